$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 244859.73
$ws.Range("J17").Value = 244859.73
$ws.Range("L17").Value = 734579.1900000001
$ws.Range("N17").Value = -734915.1900000001
$ws.Range("H26").Value = 9933.333000000001
$ws.Range("J26").Value = 9933.333000000001
$ws.Range("L26").Value = 9933.333000000001
$ws.Range("N26").Value = -10621.333
$ws.Range("H32").Value = 999.6667
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 999.5
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 999.5
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -1651.5
$ws.Range("H38").Value = 260.2
$ws.Range("I38").Value = 260.2
$ws.Range("K38").Value = 780.5999999999999
$ws.Range("M38").Value = -408.5999999999999
$ws.Range("H40").Value = 2054.3704
$ws.Range("I40").Value = 2533.75
$ws.Range("J40").Value = 1670.8667
$ws.Range("K40").Value = 2533.75
$ws.Range("L40").Value = 1670.8667
$ws.Range("M40").Value = -2358.75
$ws.Range("N40").Value = -2020.8667
$ws.Range("H58").Value = 2505
$ws.Range("J58").Value = 4972.5
$ws.Range("L58").Value = 14917.5
$ws.Range("N58").Value = -15217.5
$ws.Range("H86").Value = 3930.1365
$ws.Range("I86").Value = 3164.9167
$ws.Range("K86").Value = 3164.9167
$ws.Range("M86").Value = -2041.9167
$ws.Range("H89").Value = 3930.1365
$ws.Range("I89").Value = 3164.9167
$ws.Range("K89").Value = 15824.5835
$ws.Range("M89").Value = -10208.5835
$ws.Range("H112").Value = 1566
$ws.Range("J112").Value = 1610.8462
$ws.Range("L112").Value = 4832.5386
$ws.Range("N112").Value = -7048.5386
$ws.Range("H128").Value = 76900
$ws.Range("J128").Value = 76900
$ws.Range("L128").Value = 76900
$ws.Range("N128").Value = -86860
$ws.Range("H137").Value = 2655.3635
$ws.Range("I137").Value = 1705
$ws.Range("J137").Value = 4608.8887
$ws.Range("K137").Value = 5115
$ws.Range("L137").Value = 13826.6661
$ws.Range("M137").Value = -2565
$ws.Range("N137").Value = -18926.6661

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 5000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -4598
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 5972.091
$ws.Range("I31").Value = 5972.091
$ws.Range("K31").Value = 5972.091
$ws.Range("M31").Value = -5678.091
$ws.Range("H32").Value = 7548.269
$ws.Range("I32").Value = 6827.7607
$ws.Range("J32").Value = 14856.286
$ws.Range("K32").Value = 6827.7607
$ws.Range("L32").Value = 14856.286
$ws.Range("M32").Value = -6540.7607
$ws.Range("N32").Value = -15430.286
$ws.Range("H102").Value = 3104.0833
$ws.Range("I102").Value = 2737
$ws.Range("K102").Value = 2737
$ws.Range("M102").Value = -1115
$ws.Range("H122").Value = 4631134
$ws.Range("I122").Value = 1505.0952
$ws.Range("K122").Value = 4515.2856
$ws.Range("M122").Value = -2065.2856

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 280.0476
$ws.Range("I22").Value = 243.1875
$ws.Range("J22").Value = 398
$ws.Range("K22").Value = 243.1875
$ws.Range("L22").Value = 398
$ws.Range("M22").Value = 106.8125
$ws.Range("N22").Value = -1098
$ws.Range("H31").Value = 10927.134
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10927.134
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10927.134
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -11517.134
$ws.Range("H34").Value = 10927.134
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10927.134
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10927.134
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -11331.134
$ws.Range("H107").Value = 1151.1333
$ws.Range("I107").Value = 1162.8334
$ws.Range("J107").Value = 1104.3334
$ws.Range("K107").Value = 1162.8334
$ws.Range("L107").Value = 1104.3334
$ws.Range("M107").Value = 757.1666
$ws.Range("N107").Value = -4944.3334
$ws.Range("H132").Value = 2990.6667
$ws.Range("I132").Value = 2440.9473
$ws.Range("J132").Value = 5079.6
$ws.Range("K132").Value = 7322.841899999999
$ws.Range("L132").Value = 15238.8
$ws.Range("M132").Value = -4792.841899999999
$ws.Range("N132").Value = -20298.8
$ws.Range("H134").Value = 2967.484
$ws.Range("I134").Value = 1738.1143
$ws.Range("K134").Value = 5214.3429
$ws.Range("M134").Value = -2679.3429

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 4427.9287
$ws.Range("H30").Value = 4427.9287
$ws.Range("H44").Value = 14940
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H99").Value = 1829.2222
$ws.Range("I99").Value = 919
$ws.Range("K99").Value = 2757
$ws.Range("M99").Value = -511
$ws.Range("H113").Value = 804.1277
$ws.Range("I113").Value = 805.5278
$ws.Range("J113").Value = 799.5454999999999
$ws.Range("K113").Value = 2416.5834
$ws.Range("L113").Value = 2398.6365
$ws.Range("M113").Value = -246.5834
$ws.Range("N113").Value = -6738.6365
$ws.Range("H131").Value = 485.70706
$ws.Range("I131").Value = 282.06668
$ws.Range("J131").Value = 799
$ws.Range("K131").Value = 846.2000400000001
$ws.Range("L131").Value = 2397
$ws.Range("M131").Value = 4193.79996
$ws.Range("N131").Value = -12477

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1537.2106
$ws.Range("I97").Value = 1631.1765
$ws.Range("J97").Value = 738.5
$ws.Range("K97").Value = 1631.1765
$ws.Range("L97").Value = 738.5
$ws.Range("M97").Value = -1135.1765
$ws.Range("N97").Value = -1730.5
$ws.Range("H116").Value = 79433.336
$ws.Range("J116").Value = 79433.336
$ws.Range("L116").Value = 79433.336
$ws.Range("N116").Value = -88611.336
$ws.Range("H122").Value = 3739.88
$ws.Range("I122").Value = 4406.8125
$ws.Range("J122").Value = 2554.2222
$ws.Range("K122").Value = 13220.4375
$ws.Range("L122").Value = 7662.6666
$ws.Range("M122").Value = -10770.4375
$ws.Range("N122").Value = -12562.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5320.143
$ws.Range("I100").Value = 4275
$ws.Range("J100").Value = 7201.4
$ws.Range("K100").Value = 4275
$ws.Range("L100").Value = 7201.4
$ws.Range("M100").Value = -3734
$ws.Range("N100").Value = -8283.4
$ws.Range("H122").Value = 5222.378
$ws.Range("I122").Value = 4158.8823
$ws.Range("J122").Value = 8509.546
$ws.Range("K122").Value = 12476.6469
$ws.Range("L122").Value = 25528.638
$ws.Range("M122").Value = -10026.6469
$ws.Range("N122").Value = -30428.638

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1419.2632
$ws.Range("I107").Value = 801.8461
$ws.Range("J107").Value = 2757
$ws.Range("K107").Value = 2405.5383
$ws.Range("L107").Value = 8271
$ws.Range("M107").Value = -485.5383000000002
$ws.Range("N107").Value = -12111
$ws.Range("H122").Value = 3549.6
$ws.Range("I122").Value = 2071.1724
$ws.Range("J122").Value = 7447.273
$ws.Range("K122").Value = 6213.5172
$ws.Range("L122").Value = 22341.819
$ws.Range("M122").Value = -3763.5172
$ws.Range("N122").Value = -27241.819
$ws.Range("H123").Value = 43333
$ws.Range("I123").Value = 43333
$ws.Range("K123").Value = 43333
$ws.Range("M123").Value = -38433
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
